$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking text in column D (Price) and E (Volume) as Text format
# so values like "288.12" or "-0.42%" are stored as text, not numbers/percentages.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

# Apply cell value updates
$ws.Range("D2").Value = "288.12"
$ws.Range("E2").Value = "-0.42%"
$ws.Range("E3").Value = "2.73%"
$ws.Range("D4").Value = "4.960"
$ws.Range("E4").Value = "1.10%"
$ws.Range("D5").Value = "0.07334"
$ws.Range("E5").Value = "1.12%"
$ws.Range("D6").Value = "2.307"
$ws.Range("E6").Value = "29.20%"
$ws.Range("D7").Value = "7.713"
$ws.Range("E7").Value = "1.85%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9119"
$ws.Range("E8").Value = "1.19%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.09164"
$ws.Range("E9").Value = "17.18%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1699"
$ws.Range("E10").Value = "2.00%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.08238"
$ws.Range("E11").Value = "3.64%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03112"
$ws.Range("E12").Value = "2.28%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09974"
$ws.Range("E13").Value = "-0.42%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001497"
$ws.Range("E14").Value = "-0.09%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005717"
$ws.Range("E15").Value = "-0.87%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.472"
$ws.Range("E16").Value = "0.03%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "3.741"
$ws.Range("E17").Value = "0.75%"
$ws.Range("D18").Value = "2.039"
$ws.Range("E18").Value = "-1.86%"
$ws.Range("E19").Value = "0.08%"
$ws.Range("E20").Value = "-0.51%"
$ws.Range("D21").Value = "4.166"
$ws.Range("E21").Value = "5.12%"
$ws.Range("D22").Value = "0.2124"
$ws.Range("E22").Value = "-2.48%"
$ws.Range("D23").Value = "0.04520"
$ws.Range("E23").Value = "0.47%"
$ws.Range("E24").Value = "-0.24%"
$ws.Range("D25").Value = "0.004174"
$ws.Range("E25").Value = "-5.64%"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").Value = "0.16%"
$ws.Range("E27").Value = "0.07%"
$ws.Range("D39").Value = "0.01568"
$ws.Range("E39").Value = "0.42%"
$ws.Range("D40").Value = "0.04464"
$ws.Range("E40").Value = "3.13%"
$ws.Range("D41").Value = "0.007346"
$ws.Range("E41").Value = "-0.26%"
$ws.Range("D42").Value = "0.009862"
$ws.Range("E42").Value = "-2.06%"
$ws.Range("D43").Value = "0.1330"
$ws.Range("E43").Value = "1.96%"
$ws.Range("D44").Value = "0.002241"
$ws.Range("E44").Value = "11.87%"
$ws.Range("D45").Value = "0.008748"
$ws.Range("E45").Value = "-6.28%"
$ws.Range("D46").Value = "0.00006121"
$ws.Range("E46").Value = "3.94%"
$ws.Range("E47").Value = "0.08%"
$ws.Range("D48").Value = "2.445"
$ws.Range("E48").Value = "8.46%"
$ws.Range("D49").Value = "0.002000"
$ws.Range("E49").Value = "-30.99%"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "0.08%"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "0.08%"

# Reset style so no stray formatting/quote-prefix style remains on the range
$priceVolRange.Style = "Normal"

Write-Host "Applied all cell updates"
